# Refresh item-price-derived profit figures across the Leve profit sheets.
# (Mirrors a scheduled market-data re-pull: only the H/I/J/K/L/M/N price &
#  profit columns move; no structural changes.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I12").Value = 147.42857
$ws.Range("K12").Value = 147.42857
$ws.Range("H12").Value = 154
$ws.Range("M12").Value = 22.57142999999999
$ws.Range("I33").Value = 20834246
$ws.Range("K33").Value = 20834246
$ws.Range("H33").Value = 14708039
$ws.Range("M33").Value = -20834017
$ws.Range("I96").Value = 1949.5555
$ws.Range("K96").Value = 5848.666499999999
$ws.Range("N96").Value = -11757.4
$ws.Range("H96").Value = 2326.0715
$ws.Range("J96").Value = 3003.8
$ws.Range("L96").Value = 9011.400000000001
$ws.Range("M96").Value = -4475.666499999999
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("N100").Value = -3077
$ws.Range("H100").Value = 1663.3334
$ws.Range("J100").Value = 1995
$ws.Range("L100").Value = 1995
$ws.Range("M100").Value = -459
$ws.Range("I135").Value = 2378.2
$ws.Range("K135").Value = 21403.8
$ws.Range("N135").Value = -10069.5
$ws.Range("H135").Value = 1384
$ws.Range("J135").Value = 555.5
$ws.Range("L135").Value = 4999.5
$ws.Range("M135").Value = -18868.8
$ws.Range("I137").Value = 1283.8823
$ws.Range("K137").Value = 3851.6469
$ws.Range("H137").Value = 9178.963
$ws.Range("M137").Value = -1301.6469
$ws.Range("I138").Value = 1138.1
$ws.Range("K138").Value = 3414.3
$ws.Range("H138").Value = 3564.3103
$ws.Range("M138").Value = 1725.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I61").Value = 3506.2917
$ws.Range("K61").Value = 3506.2917
$ws.Range("H61").Value = 77483.19500000001
$ws.Range("M61").Value = -3294.2917
$ws.Range("I74").Value = 1219.7916
$ws.Range("K74").Value = 1219.7916
$ws.Range("N74").Value = -22609.5
$ws.Range("H74").Value = 9076.475
$ws.Range("J74").Value = 20861.5
$ws.Range("L74").Value = 20861.5
$ws.Range("M74").Value = -345.7916
$ws.Range("I77").Value = 1219.7916
$ws.Range("K77").Value = 6098.958000000001
$ws.Range("N77").Value = -113043.5
$ws.Range("H77").Value = 9076.475
$ws.Range("J77").Value = 20861.5
$ws.Range("L77").Value = 104307.5
$ws.Range("M77").Value = -1730.958000000001
$ws.Range("I97").Value = 2330.5
$ws.Range("K97").Value = 2330.5
$ws.Range("H97").Value = 2281.7144
$ws.Range("M97").Value = -1834.5
$ws.Range("N101").Value = -56487.5
$ws.Range("H101").Value = 49997.5
$ws.Range("J101").Value = 49997.5
$ws.Range("L101").Value = 49997.5
$ws.Range("N132").Value = -60129878
$ws.Range("H132").Value = 6682366.5
$ws.Range("J132").Value = 20041606
$ws.Range("L132").Value = 60124818
$ws.Range("N135").Value = -160738.14
$ws.Range("H135").Value = 150598.14
$ws.Range("J135").Value = 150598.14
$ws.Range("L135").Value = 150598.14
$ws.Range("I136").Value = 3506.2917
$ws.Range("K136").Value = 10518.8751
$ws.Range("H136").Value = 77483.19500000001
$ws.Range("M136").Value = -7968.875100000001
$ws.Range("N139").Value = -95279.57000000001
$ws.Range("H139").Value = 84999.57000000001
$ws.Range("J139").Value = 84999.57000000001
$ws.Range("L139").Value = 84999.57000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N94").Value = -1516.8
$ws.Range("H94").Value = 1030.3667
$ws.Range("J94").Value = 614.8
$ws.Range("L94").Value = 614.8
$ws.Range("I99").Value = 72670
$ws.Range("K99").Value = 72670
$ws.Range("N99").Value = -33245
$ws.Range("H99").Value = 48429.43
$ws.Range("J99").Value = 30249
$ws.Range("L99").Value = 30249
$ws.Range("M99").Value = -71172
$ws.Range("I105").Value = 1337.0476
$ws.Range("K105").Value = 1337.0476
$ws.Range("H105").Value = 1196.1428
$ws.Range("M105").Value = 409.9523999999999
$ws.Range("I134").Value = 28068.775
$ws.Range("K134").Value = 84206.32500000001
$ws.Range("H134").Value = 27495.46
$ws.Range("M134").Value = -81671.32500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I7").Value = 479.125
$ws.Range("K7").Value = 479.125
$ws.Range("N7").Value = -420.16667
$ws.Range("H7").Value = 422.13333
$ws.Range("J7").Value = 194.16667
$ws.Range("L7").Value = 194.16667
$ws.Range("M7").Value = -366.125
$ws.Range("I31").Value = 1621.8
$ws.Range("K31").Value = 1621.8
$ws.Range("H31").Value = 13172.154
$ws.Range("M31").Value = -1326.8
$ws.Range("I34").Value = 1621.8
$ws.Range("K34").Value = 1621.8
$ws.Range("H34").Value = 13172.154
$ws.Range("M34").Value = -1419.8
$ws.Range("N53").Value = -54126.266
$ws.Range("H53").Value = 52912.266
$ws.Range("J53").Value = 52912.266
$ws.Range("L53").Value = 52912.266
$ws.Range("I105").Value = 25464.75
$ws.Range("K105").Value = 25464.75
$ws.Range("H105").Value = 12075.8
$ws.Range("M105").Value = -23717.75
$ws.Range("I132").Value = 3852.889
$ws.Range("K132").Value = 11558.667
$ws.Range("H132").Value = 83344824
$ws.Range("M132").Value = -9028.667000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I131").Value = 12223.5
$ws.Range("K131").Value = 36670.5
$ws.Range("N131").Value = -43318.179
$ws.Range("H131").Value = 11178.881
$ws.Range("J131").Value = 11079.393
$ws.Range("L131").Value = 33238.179
$ws.Range("M131").Value = -31630.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N70").Value = -5085.1
$ws.Range("H70").Value = 4623.609
$ws.Range("J70").Value = 4545.1
$ws.Range("L70").Value = 4545.1
$ws.Range("N73").Value = -6417.1
$ws.Range("H73").Value = 4623.609
$ws.Range("J73").Value = 4545.1
$ws.Range("L73").Value = 4545.1
$ws.Range("I107").Value = 499.5
$ws.Range("K107").Value = 499.5
$ws.Range("H107").Value = 675
$ws.Range("M107").Value = 1420.5
$ws.Range("I113").Value = 1875
$ws.Range("K113").Value = 1875
$ws.Range("H113").Value = 1996.5333
$ws.Range("M113").Value = 295
$ws.Range("I126").Value = 7380.909
$ws.Range("K126").Value = 22142.727
$ws.Range("H126").Value = 7410.769
$ws.Range("M126").Value = -19672.727
$ws.Range("I132").Value = 1707.25
$ws.Range("K132").Value = 5121.75
$ws.Range("H132").Value = 9055.842000000001
$ws.Range("M132").Value = -2591.75
$ws.Range("N140").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N141").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I100").Value = 2683.6155
$ws.Range("K100").Value = 2683.6155
$ws.Range("N100").Value = -4272.5557
$ws.Range("H100").Value = 2891
$ws.Range("J100").Value = 3190.5557
$ws.Range("L100").Value = 3190.5557
$ws.Range("M100").Value = -2142.6155
$ws.Range("I132").Value = 5279.5713
$ws.Range("K132").Value = 15838.7139
$ws.Range("N132").Value = -9276691.399999999
$ws.Range("H132").Value = 2010701.2
$ws.Range("J132").Value = 3090543.8
$ws.Range("L132").Value = 9271631.399999999
$ws.Range("M132").Value = -13308.7139
$ws.Range("I136").Value = 31274.285
$ws.Range("K136").Value = 93822.855
$ws.Range("H136").Value = 232442.4
$ws.Range("M136").Value = -91272.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N46").Value = -132176.5
$ws.Range("H46").Value = 131714.5
$ws.Range("J46").Value = 131714.5
$ws.Range("L46").Value = 131714.5
$ws.Range("I132").Value = 4135.968
$ws.Range("K132").Value = 12407.904
$ws.Range("H132").Value = 9336.5
$ws.Range("M132").Value = -9877.903999999999
$ws.Range("N134").Value = -400213.5
$ws.Range("H134").Value = 131714.5
$ws.Range("J134").Value = 131714.5
$ws.Range("L134").Value = 395143.5
$ws.Range("I136").Value = 1439.138
$ws.Range("K136").Value = 4317.414
$ws.Range("H136").Value = 11539.706
$ws.Range("M136").Value = -1767.414
